$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AdWatchGlobalConfig")

# Update data values
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 20
$ws.Range("D5").Value = 21

# Update selected cell / view
$ws.Range("E18").Select()
